# Add files via upload
#
# The sheet "caña" (first sheet) already has an IMAGEN column (E) populated
# for rows 2, 4, 6-22 with base64 data-URI "pictures" stored as plain text.
# Two more products (rows 3 and 5) need the same picture that row 2 already
# carries, and all of the remaining rows (23-111) need the same picture that
# row 22 already carries, in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Fill in the two previously-blank IMAGEN cells with the same image as E2 ---
$img1 = $ws.Range("E2").Value()
$ws.Range("E3").Value = $img1
$ws.Range("E5").Value = $img1

# --- Fill in IMAGEN for every remaining data row (23-111) with the same image as E22 ---
$img2 = $ws.Range("E22").Value()
for ($r = 23; $r -le 111; $r++) {
    $ws.Cells.Item($r, 5).Value = $img2
}

# --- Update the view/selection on the sheet to match where the edit left off ---
$ws.Activate()
$ws.Range("E22:E111").Select()
